# The document contains <id>...</id> markers that were each split across
# three separate runs (the literal "<id>" tag, the page id text, and the
# literal "</id>" tag) because the page id text carried different run
# formatting from the surrounding Courier-New/olive-colored tag runs.
#
# This edit "downloads"/normalizes those markers into a single run per
# marker, merging the three runs' text into one contiguous string while
# keeping the formatting of the first ("<id>") run - exactly what Word's
# Find & Replace does when a replacement spans a run boundary: the
# replaced text collapses into one run carrying the start run's rPr.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>p135v_1</id>", $false, $false, $false, $false, $false,
    $true, 1, $false, "<id>p135v_1</id>", 2) | Out-Null

$d.Content.Find.Execute(
    "<id>p136r_1</id>", $false, $false, $false, $false, $false,
    $true, 1, $false, "<id>p136r_1</id>", 2) | Out-Null
